$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.67670000000003
$ws.Range("E4").Value = 13.9729

$ws.Range("E5").Value = 12.50749999999999

$ws.Range("A6").Value = -21.49240000000002

$ws.Range("A7").Value = -21.50290000000001

$ws.Range("E8").Value = 14.2111

$ws.Range("A16").Value = -20.12959999999999
$ws.Range("E16").Value = 11.8808

$ws.Range("A20").Value = -22.11540000000003

$ws.Range("E22").Value = 11.3752
